$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The docente-name-only row (no label in column A) is removed; everything
# below it shifts up by one row.
$ws.Rows.Item(13).Delete()

# After the row shift, a handful of fields were re-populated with new
# (mostly short) values for the B/C columns.
$ws.Range("B10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
